$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "Medical records"
$ws.Range("B2").Value = "../medical-records.pdf"
$ws.Range("C2").Value = "pdf"

$ws.Range("F5").Value = "PMAS5 Appendix A6"
$ws.Range("F6").Value = "PMAS5 Ch 8.3"
$ws.Range("F3").Value = "PMAS5 Appendix A3"

$ws.Range("D9").Select()
